$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-12-13", "overview", "K02000001", "United Kingdom", 10873468, 54661, 38, 146477),
    @("2021-12-14", "overview", "K02000001", "United Kingdom", 10932545, 59610, 150, 146627),
    @("2021-12-15", "overview", "K02000001", "United Kingdom", 11010286, 78610, 165, 146791),
    @("2021-12-16", "overview", "K02000001", "United Kingdom", 11097851, 88376, 146, 146937),
    @("2021-12-17", "overview", "K02000001", "United Kingdom", 11190354, 93045, 111, 147048),
    @("2021-12-18", "overview", "K02000001", "United Kingdom", 11279428, 90418, 125, 147173),
    @("2021-12-19", "overview", "K02000001", "United Kingdom", 11361387, 82886, 45, 147218)
)

$startRow = 489
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    # Leading apostrophe forces the date-shaped string to be kept as literal
    # text instead of Excel auto-converting it to a date serial number.
    $ws.Cells.Item($row, 1).Value = "'" + $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}

# Clear the implicit "quote prefix" cell style that typing a leading
# apostrophe applies, so column A ends up with the default (unstyled) cell
# format, matching plain text cells elsewhere in the sheet.
$ws.Range("A$startRow`:A$endRow").Style = "Normal"
